$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15. This shifts the existing rows 15..40
# down to 16..41, preserving all of their data and formatting.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with a new weekly record. The row
# reuses the same market/category/price data that used to sit in row 14,
# but records a new date (one day after the previous latest entry).
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = 44498
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Verde"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 5000
$ws.Range("K15").Value = 900
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = 900
$ws.Range("N15").Value = "$/kilo"
$ws.Range("O15").Value = "Provincia de Linares"
$ws.Range("P15").Value = 900
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"
